$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update column widths (stored OOXML width = ColumnWidth + 5/6,
# so subtract the 5/6 padding offset to land on the exact target stored widths)
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws.Columns.Item(5).ColumnWidth = 22.166666666666668
$ws.Columns.Item(6).ColumnWidth = 22.166666666666668

# Row 3 updates
$ws.Range("D3").Value = 12772.25
$ws.Range("E3").Value = 951.0900000000001
$ws.Range("F3").Value = 0.930695442946105

# Row 4 updates
$ws.Range("D4").Value = 18171.67
$ws.Range("E4").Value = -4448.33
$ws.Range("F4").Value = 1.324143393663642
